$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 had previously only blank formatted cells in C22:D22. The new
# review row reuses the same look-and-feel as the row above it (row 21),
# so copy that row's formatting down before filling in the values.
$ws.Range("A21:G21").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)

$ws.Range("A22").Value = "com.hamxa.shaynachim"
$ws.Range("B22").Value = "bitcoin"
$ws.Range("C22").Value = "sinuspai@gmail.com"
$ws.Range("D22").Value = "armonravid2@gmail.com "
$ws.Range("E22").Value = "27/5/2019 15:59"
$ws.Range("F22").Value = "welcome to the bitcoin jungle. Now make money"
$ws.Range("G22").Value = "no"

# New review's email gets the same mailto hyperlink treatment as the rest
# of the email column.
$ws.Hyperlinks.Add($ws.Range("C22"), "mailto:sinuspai@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "sinuspai@gmail.com")

$ws.Range("G23").Select()
